# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.940.22'
$ws.Range('E2').Value = '  +1.50%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.748.70'
$ws.Range('E3').Value = '  -0.81%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '334.44'
$ws.Range('E5').Value = '  -0.39%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.07%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3849'
$ws.Range('E7').Value = '  +0.30%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.99'
$ws.Range('E9').Value = '  -1.95%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.112'
$ws.Range('E10').Value = '  -2.02%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07189'
$ws.Range('E11').Value = '  -2.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.38'
$ws.Range('E13').Value = '  +0.08%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.143'
$ws.Range('E14').Value = '  -2.97%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.750.17'
$ws.Range('E15').Value = '  -0.80%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.077'
$ws.Range('E16').Value = '  +0.29%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001056'
$ws.Range('E17').Value = '  -1.30%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06612'
$ws.Range('E18').Value = '  -0.79%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '78.99'
$ws.Range('E19').Value = '  -3.65%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.13%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.72'
$ws.Range('E21').Value = '  -3.24%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.168'
$ws.Range('E22').Value = '  -3.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.933.18'
$ws.Range('E23').Value = '  +1.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.61'
$ws.Range('E24').Value = '  -3.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.404'
$ws.Range('E25').Value = '  +0.60%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.18'
$ws.Range('E26').Value = '  +0.75%  '

$ws.Range('E27').Value = '  -3.80%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.292'
$ws.Range('E28').Value = '  -4.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.950.38'
$ws.Range('E29').Value = '  -0.72%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.281'
$ws.Range('E30').Value = '  -9.85%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '130.56'
$ws.Range('E31').Value = '  -2.87%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.026'
$ws.Range('E32').Value = '  +1.63%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.801'
$ws.Range('E33').Value = '  -4.55%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08785'
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.13'
$ws.Range('E35').Value = '  -4.32%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.538'
$ws.Range('E36').Value = '  +1.82%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6508'
$ws.Range('E37').Value = '  -3.80%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02277'
$ws.Range('E38').Value = '  -5.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.128'
$ws.Range('E39').Value = '  -3.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06114'
$ws.Range('E40').Value = '  -2.96%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2089'
$ws.Range('E41').Value = '  -3.91%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.203'
$ws.Range('E42').Value = '  -2.88%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.946'
$ws.Range('E43').Value = '  -3.35%  '

$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.08%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.63'
$ws.Range('E45').Value = '  -3.40%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.829'
$ws.Range('E46').Value = '  +0.16%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6007'
$ws.Range('E47').Value = '  -3.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.69'
$ws.Range('E48').Value = '  -2.80%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.991'
$ws.Range('E49').Value = '  -3.67%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.165'
$ws.Range('E50').Value = '  +1.85%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.108'
$ws.Range('E51').Value = '  +4.64%  '
